# Insert a new weekly price record as row 139 of the "Jengibre" data sheet.
# This pushes the previously existing rows 139-185 down to 140-186,
# growing the used range from A1:R185 to A1:R186.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 139..185 down one row, duplicating row 139's formatting
# (including the date-style cell D139) into the freshly inserted row.
$ws.Rows.Item(139).Insert()

# Populate the new row 139 with the new weekly observation.
$ws.Cells.Item(139, 1).Value  = 10
$ws.Cells.Item(139, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(139, 3).Value  = "La Araucanía"
$ws.Cells.Item(139, 4).Value  = 44795
$ws.Cells.Item(139, 5).Value  = 9
$ws.Cells.Item(139, 6).Value  = 100114007
$ws.Cells.Item(139, 7).Value  = "Jengibre"
$ws.Cells.Item(139, 8).Value  = "Sin especificar"
$ws.Cells.Item(139, 9).Value  = "Primera"
$ws.Cells.Item(139, 10).Value = 100
$ws.Cells.Item(139, 11).Value = 20000
$ws.Cells.Item(139, 12).Value = 20000
$ws.Cells.Item(139, 13).Value = 20000
$ws.Cells.Item(139, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(139, 15).Value = "Perú"
$ws.Cells.Item(139, 16).Value = 1538
$ws.Cells.Item(139, 17).Value = 13
$ws.Cells.Item(139, 18).Value = "Hortaliza"
